$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 69. This shifts the existing rows 69:176
# down to 71:178, and the sheet dimension grows from A1:R176 to A1:R178.
$ws.Rows("69:70").Insert()

# Fill in the two new rows (69 and 70) with their data.

# Row 69
$ws.Cells.Item(69, 1).Value = 7
$ws.Cells.Item(69, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(69, 3).Value = "Ñuble"
$ws.Cells.Item(69, 4).Value = 44803
$ws.Cells.Item(69, 5).Value = 16
$ws.Cells.Item(69, 6).Value = 100112045
$ws.Cells.Item(69, 7).Value = "Zapallo"
$ws.Cells.Item(69, 8).Value = "Camote"
$ws.Cells.Item(69, 9).Value = "1a (guarda)"
$ws.Cells.Item(69, 10).Value = 120
$ws.Cells.Item(69, 11).Value = 800
$ws.Cells.Item(69, 12).Value = 900
$ws.Cells.Item(69, 13).Value = 850
$ws.Cells.Item(69, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(69, 15).Value = "Región del Maule"
$ws.Cells.Item(69, 16).Value = 850
$ws.Cells.Item(69, 17).Value = 1
$ws.Cells.Item(69, 18).Value = "Hortaliza"

# Row 70
$ws.Cells.Item(70, 1).Value = 7
$ws.Cells.Item(70, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(70, 3).Value = "Ñuble"
$ws.Cells.Item(70, 4).Value = 44803
$ws.Cells.Item(70, 5).Value = 16
$ws.Cells.Item(70, 6).Value = 100112045
$ws.Cells.Item(70, 7).Value = "Zapallo"
$ws.Cells.Item(70, 8).Value = "Camote"
$ws.Cells.Item(70, 9).Value = "2a (guarda)"
$ws.Cells.Item(70, 10).Value = 80
$ws.Cells.Item(70, 11).Value = 700
$ws.Cells.Item(70, 12).Value = 700
$ws.Cells.Item(70, 13).Value = 700
$ws.Cells.Item(70, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(70, 15).Value = "Región del Maule"
$ws.Cells.Item(70, 16).Value = 700
$ws.Cells.Item(70, 17).Value = 1
$ws.Cells.Item(70, 18).Value = "Hortaliza"
